$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AEDB.CEA"
$ws.Range("B2").Value = "epmajor.3years"
$ws.Range("C2").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D2").Value = 0.109110457720438
$ws.Range("E2").Value = 0.185236502181912
$ws.Range("F2").Value = 1.11528553543583
$ws.Range("G2").Value = 0.775728193273921
$ws.Range("H2").Value = 1.60347636754407
$ws.Range("I2").Value = 0.589033243638374
$ws.Range("J2").Value = 0.55583897327207
$ws.Range("K2").Value = 1186
$ws.Range("L2").Value = 139

$ws.Range("A3").Value = "AEDB.CEA"
$ws.Range("B3").Value = "epmajor.3years"
$ws.Range("C3").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D3").Value = -0.0124325654075699
$ws.Range("E3").Value = 0.189652402427621
$ws.Range("F3").Value = 0.987644399645869
$ws.Range("G3").Value = 0.681028397384233
$ws.Range("H3").Value = 1.43230658794615
$ws.Range("I3").Value = -0.0655544841427181
$ws.Range("J3").Value = 0.94773252751171
$ws.Range("K3").Value = 1187
$ws.Range("L3").Value = 140

$ws.Range("A4").Value = "AEDB.CEA"
$ws.Range("B4").Value = "epmajor.3years"
$ws.Range("C4").Value = "MCP1_rank"
$ws.Range("D4").Value = -0.27321691717627
$ws.Range("E4").Value = 0.253769136620049
$ws.Range("F4").Value = 0.760927711435141
$ws.Range("G4").Value = 0.462733296218974
$ws.Range("H4").Value = 1.25128445858783
$ws.Range("I4").Value = -1.07663572022684
$ws.Range("J4").Value = 0.281643039650307
$ws.Range("K4").Value = 549
$ws.Range("L4").Value = 70

$ws.Range("A5").Value = "AEDB.CEA"
$ws.Range("B5").Value = "epstroke.3years"
$ws.Range("C5").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D5").Value = 0.305916951796889
$ws.Range("E5").Value = 0.257882219168381
$ws.Range("F5").Value = 1.35786953324035
$ws.Range("G5").Value = 0.819113846746521
$ws.Range("H5").Value = 2.25098095536523
$ws.Range("I5").Value = 1.18626616749076
$ws.Range("J5").Value = 0.235517203065819
$ws.Range("K5").Value = 1186
$ws.Range("L5").Value = 73

$ws.Range("A6").Value = "AEDB.CEA"
$ws.Range("B6").Value = "epstroke.3years"
$ws.Range("C6").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D6").Value = 0.0802507513231576
$ws.Range("E6").Value = 0.260754086232535
$ws.Range("F6").Value = 1.08355873739989
$ws.Range("G6").Value = 0.649971178598896
$ws.Range("H6").Value = 1.80638707692635
$ws.Range("I6").Value = 0.307764117842401
$ws.Range("J6").Value = 0.758261826972502
$ws.Range("K6").Value = 1187
$ws.Range("L6").Value = 74

$ws.Range("A7").Value = "AEDB.CEA"
$ws.Range("B7").Value = "epstroke.3years"
$ws.Range("C7").Value = "MCP1_rank"
$ws.Range("D7").Value = -0.382334210661852
$ws.Range("E7").Value = 0.351409093297658
$ws.Range("F7").Value = 0.682266994193947
$ws.Range("G7").Value = 0.342632774535179
$ws.Range("H7").Value = 1.35856312052439
$ws.Range("I7").Value = -1.0880031790697
$ws.Range("J7").Value = 0.276593704098031
$ws.Range("K7").Value = 549
$ws.Range("L7").Value = 36

$ws.Range("A8").Value = "AEDB.CEA"
$ws.Range("B8").Value = "epcoronary.3years"
$ws.Range("C8").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D8").Value = -0.215530764113133
$ws.Range("E8").Value = 0.227947710052008
$ws.Range("F8").Value = 0.806113470510643
$ws.Range("G8").Value = 0.515659674924895
$ws.Range("H8").Value = 1.26017014503482
$ws.Range("I8").Value = -0.94552721790431
$ws.Range("J8").Value = 0.344389781976305
$ws.Range("K8").Value = 1186
$ws.Range("L8").Value = 91

$ws.Range("A9").Value = "AEDB.CEA"
$ws.Range("B9").Value = "epcoronary.3years"
$ws.Range("C9").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D9").Value = 0.0631563377814632
$ws.Range("E9").Value = 0.234085342666354
$ws.Range("F9").Value = 1.06519335618066
$ws.Range("G9").Value = 0.673241660804058
$ws.Range("H9").Value = 1.68533373988814
$ws.Range("I9").Value = 0.269800479868068
$ws.Range("J9").Value = 0.787313753627286
$ws.Range("K9").Value = 1187
$ws.Range("L9").Value = 91

$ws.Range("A10").Value = "AEDB.CEA"
$ws.Range("B10").Value = "epcoronary.3years"
$ws.Range("C10").Value = "MCP1_rank"
$ws.Range("D10").Value = 0.196482239218439
$ws.Range("E10").Value = 0.314183517078195
$ws.Range("F10").Value = 1.2171137037956
$ws.Range("G10").Value = 0.657495437409902
$ws.Range("H10").Value = 2.25304341852567
$ws.Range("I10").Value = 0.625374115885072
$ws.Range("J10").Value = 0.531725546344934
$ws.Range("K10").Value = 549
$ws.Range("L10").Value = 46

$ws.Range("A11").Value = "AEDB.CEA"
$ws.Range("B11").Value = "epcvdeath.3years"
$ws.Range("C11").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D11").Value = -0.029035347507777
$ws.Range("E11").Value = 0.323506909861162
$ws.Range("F11").Value = 0.971382127922432
$ws.Range("G11").Value = 0.515246995264451
$ws.Range("H11").Value = 1.83132215640154
$ws.Range("I11").Value = -0.0897518619315982
$ws.Range("J11").Value = 0.928484402343315
$ws.Range("K11").Value = 1186
$ws.Range("L11").Value = 45

$ws.Range("A12").Value = "AEDB.CEA"
$ws.Range("B12").Value = "epcvdeath.3years"
$ws.Range("C12").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D12").Value = -0.201269455840919
$ws.Range("E12").Value = 0.33297833575438
$ws.Range("F12").Value = 0.817692069962009
$ws.Range("G12").Value = 0.425748282323535
$ws.Range("H12").Value = 1.57045923386875
$ws.Range("I12").Value = -0.604452104623961
$ws.Range("J12").Value = 0.545543103119405
$ws.Range("K12").Value = 1187
$ws.Range("L12").Value = 45

$ws.Range("A13").Value = "AEDB.CEA"
$ws.Range("B13").Value = "epcvdeath.3years"
$ws.Range("C13").Value = "MCP1_rank"
$ws.Range("D13").Value = -0.0979334533145785
$ws.Range("E13").Value = 0.416021477852458
$ws.Range("F13").Value = 0.906709240240881
$ws.Range("G13").Value = 0.401183821129396
$ws.Range("H13").Value = 2.04923928393671
$ws.Range("I13").Value = -0.235404801261993
$ws.Range("J13").Value = 0.813894559905866
$ws.Range("K13").Value = 549
$ws.Range("L13").Value = 26

